# The upload re-saved the workbook with the "school list" contents cleared
# out of column A (rows 4 through 54), leaving only the first three rows'
# text in place, and reset the view back to the top of the sheet with the
# selection sitting on A3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the contents (but keep cell formatting/styles) for A4:A54 - the
# school-name rows beyond the first three get wiped out.
$ws.Range("A4:A54").ClearContents()

# Park the selection on A3, matching the saved view state (no more
# topLeftCell offset, selection on A3 instead of A37).
$ws.Range("A3").Select()
